$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table 3 data: fill in the newly-collected "H" column (rows 2-7),
# leaving the intersecting point (row 8, not yet known) for later.
$ws.Range("H2").Value = 139.369
$ws.Range("H3").Value = 246.31
$ws.Range("H4").Value = 553.84
$ws.Range("H5").Value = 1028.31
$ws.Range("H6").Value = 1977.6
$ws.Range("H7").Value = 3666.89

# Reflect where the author left the cursor / viewport after entering
# the last value: selection on H7, scrolled so column E is left-most.
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H7").Select()
